$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 4326
$ws.Range("F5").Value = 317
$ws.Range("F7").Value = 3440
$ws.Range("F8").Value = 969
$ws.Range("F11").Value = 293
$ws.Range("F12").Value = 2330
$ws.Range("F13").Value = 1250
$ws.Range("F14").Value = 28
$ws.Range("F15").Value = 1979
$ws.Range("F16").Value = 499
$ws.Range("F17").Value = 244
$ws.Range("F18").Value = 51
$ws.Range("F19").Value = 9673
$ws.Range("F20").Value = 5916
$ws.Range("F23").Value = 806
$ws.Range("F24").Value = 116
$ws.Range("F26").Value = 3503
$ws.Range("F30").Value = 101
$ws.Range("F31").Value = 227
$ws.Range("F32").Value = 204
$ws.Range("F33").Value = 4773
$ws.Range("F34").Value = 15
$ws.Range("F35").Value = 1023
$ws.Range("F36").Value = 130
$ws.Range("F37").Value = 9
$ws.Range("F38").Value = 454

$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 122
$ws.Range("F14").Value = 116
$ws.Range("F15").Value = 3513

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 8635
$ws.Range("F3").Value = 412
$ws.Range("F4").Value = 1510

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 8635
$ws.Range("F4").Value = 412
$ws.Range("F5").Value = 1510
$ws.Range("F7").Value = 4326
$ws.Range("F10").Value = 3440
$ws.Range("F11").Value = 969
$ws.Range("F14").Value = 293
$ws.Range("F15").Value = 2330
$ws.Range("F19").Value = 1250
$ws.Range("F21").Value = 28
$ws.Range("F22").Value = 122
$ws.Range("F23").Value = 499
$ws.Range("F24").Value = 244
$ws.Range("F25").Value = 51
$ws.Range("F26").Value = 9674
$ws.Range("F27").Value = 3513
$ws.Range("F31").Value = 806
$ws.Range("F32").Value = 116
$ws.Range("F34").Value = 3503
$ws.Range("F38").Value = 101
$ws.Range("F39").Value = 227
$ws.Range("F41").Value = 204
$ws.Range("F42").Value = 4773
$ws.Range("F43").Value = 1023
$ws.Range("F44").Value = 130
$ws.Range("F45").Value = 454
